# Daily Report update: 2026-01-09
# Adds the new business-day block (date serial 46030 / 2026-01-08 data)
# to Daily_Data, and refreshes the downstream Today_Summary and
# Monthly_Stats roll-up sheets to reflect the corrected BRINK'S, INC.
# "Eligible" figures that came in with this day's data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append rows 90-111 (date 46030) after the existing data
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(90, 46030, "ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @(91, 46030, "ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(92, 46030, "BRINK'S, INC. Registered", 90027.72500000001, 0, 0, 0, 0, 90027.72500000001),
    @(93, 46030, "BRINK'S, INC. Eligible", 5744.711, 0, 260.167, -260.167, 0, 5484.544),
    @(94, 46030, "CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @(95, 46030, "CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(96, 46030, "DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @(97, 46030, "DELAWARE DEPOSITORY Eligible", 18509.729, 0, 0, 0, 0, 18509.729),
    @(98, 46030, "HSBC BANK, USA Registered", 1295.223, 0, 0, 0, 0, 1295.223),
    @(99, 46030, "HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @(100, 46030, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @(101, 46030, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @(102, 46030, "JP MORGAN CHASE BANK NA Registered", 124991.729, 0, 0, 0, 0, 124991.729),
    @(103, 46030, "JP MORGAN CHASE BANK NA Eligible", 125407.673, 0, 0, 0, 0, 125407.673),
    @(104, 46030, "LOOMIS INTERNATIONAL (US) LLC Registered", 68084.33, 0, 0, 0, 0, 68084.33),
    @(105, 46030, "LOOMIS INTERNATIONAL (US) LLC Eligible", 106188.481, 0, 0, 0, 0, 106188.481),
    @(106, 46030, "MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @(107, 46030, "MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(108, 46030, "MANFRA, TORDELLA & BROOKES, LLC Registered", 54605.27, 0, 0, 0, 0, 54605.27),
    @(109, 46030, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 1068.408, 0, 0, 0, 0, 1068.408),
    @(110, 46030, "STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @(111, 46030, "STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $cellA = $daily.Cells.Item($r, 1)
    $cellA.Value = $row[1]
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $daily.Cells.Item($r, 2).Value = $row[2]
    $daily.Cells.Item($r, 3).Value = $row[3]
    $daily.Cells.Item($r, 4).Value = $row[4]
    $daily.Cells.Item($r, 5).Value = $row[5]
    $daily.Cells.Item($r, 6).Value = $row[6]
    $daily.Cells.Item($r, 7).Value = $row[7]
    $daily.Cells.Item($r, 8).Value = $row[8]
}

# ---------------------------------------------------------------------
# 2) Today_Summary: BRINK'S, INC. Eligible/Total_Stock refreshed
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("Today_Summary")
$today.Range("B3").Value = 5484.544
$today.Range("D3").Value = 95512.269

# ---------------------------------------------------------------------
# 3) Monthly_Stats: month-to-date roll-up totals + BRINK'S detail row
# ---------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")
$monthly.Range("B2").Value = 265956.889
$monthly.Range("D2").Value = 624754.525

$monthly.Range("D9").Value = 260.167
$monthly.Range("E9").Value = 5484.544
